$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 253, shifting existing rows
# 253-278 down to 254-279 (dimension grows from R278 to R279).
$ws.Rows(253).Insert()

# Populate the newly inserted row 253 with the new weekly record.
$ws.Cells(253, 1).Value = 6
$ws.Cells(253, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells(253, 3).Value = "Metropolitana"
$ws.Cells(253, 4).Value = 44946
$ws.Cells(253, 5).Value = 13
$ws.Cells(253, 6).Value = 100112001
$ws.Cells(253, 7).Value = "Berenjena"
$ws.Cells(253, 8).Value = "Sin especificar"
$ws.Cells(253, 9).Value = "Primera"
$ws.Cells(253, 10).Value = 400
$ws.Cells(253, 11).Value = 20000
$ws.Cells(253, 12).Value = 22000
$ws.Cells(253, 13).Value = 20850
$ws.Cells(253, 14).Value = "$/caja 50 unidades"
$ws.Cells(253, 15).Value = "Provincia de Huasco"
$ws.Cells(253, 16).Value = 417
$ws.Cells(253, 17).Value = 50
$ws.Cells(253, 18).Value = "Hortaliza"

# Keep the date column formatted the same way as the rest of column D.
$ws.Cells(253, 4).NumberFormat = $ws.Cells(254, 4).NumberFormat
